$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.481.06'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '3.469.12'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'586.35"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').Value = "'177.05"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = "'0.629"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.47%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '3.468.82'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('D11').Value = "'6.98"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('D13').Value = '4.069.38'
$ws.Range('E13').Value = '  -1.57%  '
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').Value = "'30.34"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '66.391.55'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '3.471.78'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').Value = "'373.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.63%  '
$ws.Range('D22').Value = "'7.68"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.24%  '
$ws.Range('D23').Value = "'73.41"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  +4.62%  '
$ws.Range('D26').Value = "'0.538"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = "'10.02"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = "'0.178"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.94%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = "'0.999"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = "'5.99"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'2.00"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'23.73"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.98%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = "'0.999"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').Value = "'7.07"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.98%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = "'1.27"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.50%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = "'1.56"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = "'161.38"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Value = "'0.888"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = "'28.34"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.23%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'1.81"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.791.05'
$ws.Range('E41').Value = '  +1.86%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = "'4.52"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = "'2.58"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = "'6.46"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = "'0.0695"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = "'25.32"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = "'342.94"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.70%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = "'40.09"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = "'0.0293"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = "'0.105"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').Value = "'1.00"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.02%  '
